# Add data for 2025-05-22
# Updates the 2025 (column L) violent-crime running totals across the
# Citywide Totals, By Neighborhood, and per-neighborhood sheets, plus a
# couple of 2023 (column J) corrections called out by the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2409
$ws.Range("L3").Value = 2429
$ws.Range("J4").Value = 1866
$ws.Range("L4").Value = 664
$ws.Range("L6").Value = 2212
$ws.Range("J7").Value = 29340
$ws.Range("L7").Value = 7859

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 142
$ws.Range("L3").Value = 166
$ws.Range("L4").Value = 38
$ws.Range("L7").Value = 498

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 58
$ws.Range("L6").Value = 43
$ws.Range("L7").Value = 191

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 97
$ws.Range("L6").Value = 125
$ws.Range("L7").Value = 359

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 85
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 288

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 47
$ws.Range("L7").Value = 145

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 124

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 65
$ws.Range("L7").Value = 253
$ws.Range("L8").Value = 498
$ws.Range("L18").Value = 58
$ws.Range("L19").Value = 223
$ws.Range("L20").Value = 200
$ws.Range("L22").Value = 25
$ws.Range("L25").Value = 45
$ws.Range("L27").Value = 81
$ws.Range("L33").Value = 359
$ws.Range("L37").Value = 288
$ws.Range("L42").Value = 251
$ws.Range("L51").Value = 90
$ws.Range("L52").Value = 155
$ws.Range("L54").Value = 160
$ws.Range("J63").Value = 218
$ws.Range("L63").Value = 27
$ws.Range("L64").Value = 52
$ws.Range("L65").Value = 145
$ws.Range("L67").Value = 291
$ws.Range("L75").Value = 31
$ws.Range("L76").Value = 91
$ws.Range("L83").Value = 191
$ws.Range("L85").Value = 415
$ws.Range("L88").Value = 107
$ws.Range("L91").Value = 113
$ws.Range("L94").Value = 93
$ws.Range("L95").Value = 111
$ws.Range("L96").Value = 78
$ws.Range("L97").Value = 70
$ws.Range("L99").Value = 124
$ws.Range("J101").Value = 29340
$ws.Range("L101").Value = 7859

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 87
$ws.Range("L3").Value = 98
$ws.Range("L4").Value = 24
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 291

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 35
$ws.Range("L6").Value = 83
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 69
$ws.Range("L7").Value = 223

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 77
$ws.Range("L6").Value = 81
$ws.Range("L7").Value = 251

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 78

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 43
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 63
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L6").Value = 72
$ws.Range("L7").Value = 253

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 19
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L3").Value = 36
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 26
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 124
$ws.Range("L3").Value = 167
$ws.Range("L6").Value = 83
$ws.Range("L7").Value = 415

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 155
